# Append new service-map rows to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# First batch (rows 10-13): ORIGINAL_SERVICE column, then SERVICE_NAME column,
# then SERVICE_ABBR column -- matches how the values were pasted in originally
# (column-at-a-time), which drives the shared-string insertion order.
$ws.Cells.Item(10, 1).Value = 113
$ws.Cells.Item(11, 1).Value = 114
$ws.Cells.Item(12, 1).Value = 51
$ws.Cells.Item(13, 1).Value = 62

$ws.Cells.Item(10, 2).Value = "USPS Priority Mail Standard"
$ws.Cells.Item(11, 2).Value = "USPS Express Mail Standard"
$ws.Cells.Item(12, 2).Value = "USPS Priority Mail"
$ws.Cells.Item(13, 2).Value = "USPS Express Mail"

$ws.Cells.Item(10, 3).Value = "PMIST"
$ws.Cells.Item(11, 3).Value = "EMIST"
$ws.Cells.Item(12, 3).Value = "PMIPR"
$ws.Cells.Item(13, 3).Value = "PMEI"

# Second batch (rows 14-17): entered row by row (A, B, C each row).
$rows = @(
    @(109, "ETOE Unregistered",   "ETUR"),
    @(110, "ETOE Registered",     "ETR"),
    @(111, "ETOE Packet Tracked", "ETPT"),
    @(112, "ETOE Parcel",         "ETP")
)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
}

# Clear the lingering selection anchor on B5 left over from authoring, and
# reselect A1 so the saved sheetView has no stale <selection> element.
$ws.Range("A1").Select()
